# Auto-generated edit script: applies Phantom_Profits.xlsx numeric updates
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2422.6365
$ws.Range("I28").Value = 2301.4666
$ws.Range("J28").Value = 2682.2856
$ws.Range("K28").Value = 2301.4666
$ws.Range("L28").Value = 2682.2856
$ws.Range("M28").Value = -1816.4666
$ws.Range("N28").Value = -3652.2856
$ws.Range("H87").Value = 157500
$ws.Range("J87").Value = 157500
$ws.Range("L87").Value = 157500
$ws.Range("N87").Value = -159996
$ws.Range("H90").Value = 157500
$ws.Range("J90").Value = 157500
$ws.Range("L90").Value = 472500
$ws.Range("N90").Value = -484980
$ws.Range("H116").Value = 7666.1113
$ws.Range("I116").Value = 7499.5
$ws.Range("K116").Value = 7499.5
$ws.Range("M116").Value = -4057.5
$ws.Range("H125").Value = 1163.3334
$ws.Range("J125").Value = 1163.3334
$ws.Range("L125").Value = 10470.0006
$ws.Range("N125").Value = -15390.0006

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5605.2
$ws.Range("I45").Value = 5503
$ws.Range("K45").Value = 5503
$ws.Range("M45").Value = -5126
$ws.Range("H122").Value = 1609
$ws.Range("I122").Value = 1609
$ws.Range("K122").Value = 4827
$ws.Range("M122").Value = -2377
$ws.Range("H124").Value = 72809.336
$ws.Range("J124").Value = 72809.336
$ws.Range("L124").Value = 72809.336
$ws.Range("N124").Value = -82629.336
$ws.Range("H135").Value = 39829.25
$ws.Range("J135").Value = 39829.25
$ws.Range("L135").Value = 39829.25
$ws.Range("N135").Value = -49969.25
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 69999
$ws.Range("J35").Value = 69999
$ws.Range("L35").Value = 69999
$ws.Range("N35").Value = -70619
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H80").Value = 1356.909
$ws.Range("I80").Value = 592.8
$ws.Range("J80").Value = 8998
$ws.Range("K80").Value = 592.8
$ws.Range("L80").Value = 8998
$ws.Range("M80").Value = 405.2
$ws.Range("N80").Value = -10994
$ws.Range("H83").Value = 1356.909
$ws.Range("I83").Value = 592.8
$ws.Range("J83").Value = 8998
$ws.Range("K83").Value = 2964
$ws.Range("L83").Value = 44990
$ws.Range("M83").Value = 2028
$ws.Range("N83").Value = -54974
$ws.Range("H86").Value = 22701.834
$ws.Range("I86").Value = 29608.5
$ws.Range("J86").Value = 8888.5
$ws.Range("K86").Value = 29608.5
$ws.Range("L86").Value = 8888.5
$ws.Range("M86").Value = -28485.5
$ws.Range("N86").Value = -11134.5
$ws.Range("H89").Value = 22701.834
$ws.Range("I89").Value = 29608.5
$ws.Range("J89").Value = 8888.5
$ws.Range("K89").Value = 148042.5
$ws.Range("L89").Value = 44442.5
$ws.Range("M89").Value = -142426.5
$ws.Range("N89").Value = -55674.5
$ws.Range("H134").Value = 2220.5
$ws.Range("I134").Value = 2220.5
$ws.Range("K134").Value = 6661.5
$ws.Range("M134").Value = -4126.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 3800
$ws.Range("I38").Value = 4375
$ws.Range("J38").Value = 1500
$ws.Range("K38").Value = 4375
$ws.Range("L38").Value = 1500
$ws.Range("M38").Value = -3998
$ws.Range("N38").Value = -2254
$ws.Range("H46").Value = 3800
$ws.Range("I46").Value = 4375
$ws.Range("J46").Value = 1500
$ws.Range("K46").Value = 4375
$ws.Range("L46").Value = 1500
$ws.Range("M46").Value = -4164
$ws.Range("N46").Value = -1922
$ws.Range("H58").Value = 3070.6924
$ws.Range("I58").Value = 2212.3
$ws.Range("K58").Value = 2212.3
$ws.Range("M58").Value = -2009.3
$ws.Range("H86").Value = 3921.3
$ws.Range("I86").Value = 3929.125
$ws.Range("K86").Value = 3929.125
$ws.Range("M86").Value = -2806.125
$ws.Range("H89").Value = 3921.3
$ws.Range("I89").Value = 3929.125
$ws.Range("K89").Value = 19645.625
$ws.Range("M89").Value = -14029.625
$ws.Range("H107").Value = 2619.8
$ws.Range("I107").Value = 2619.8
$ws.Range("K107").Value = 2619.8
$ws.Range("M107").Value = -699.8000000000002
$ws.Range("H136").Value = 3070.6924
$ws.Range("I136").Value = 2212.3
$ws.Range("K136").Value = 6636.900000000001
$ws.Range("M136").Value = -4086.900000000001
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 79
$ws.Range("I2").Value = 104.44444
$ws.Range("J2").Value = 50.375
$ws.Range("K2").Value = 626.66664
$ws.Range("L2").Value = 302.25
$ws.Range("M2").Value = -513.66664
$ws.Range("N2").Value = -528.25
$ws.Range("H39").Value = 5394.3335
$ws.Range("J39").Value = 5793.909
$ws.Range("L39").Value = 17381.727
$ws.Range("N39").Value = -17969.727
$ws.Range("H122").Value = 1394.3334
$ws.Range("J122").Value = 1785
$ws.Range("L122").Value = 16065
$ws.Range("N122").Value = -20965
$ws.Range("H130").Value = 1845
$ws.Range("I130").Value = 1690
$ws.Range("J130").Value = 2000
$ws.Range("K130").Value = 5070
$ws.Range("L130").Value = 6000
$ws.Range("M130").Value = -50
$ws.Range("N130").Value = -16040
$ws.Range("H131").Value = 825
$ws.Range("I131").Value = 825
$ws.Range("K131").Value = 2475
$ws.Range("M131").Value = 2565

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 12000
$ws.Range("J44").Value = 12000
$ws.Range("L44").Value = 12000
$ws.Range("N44").Value = -13192
$ws.Range("H113").Value = 2992.8333
$ws.Range("I113").Value = 2993.4
$ws.Range("J113").Value = 2990
$ws.Range("K113").Value = 2993.4
$ws.Range("L113").Value = 2990
$ws.Range("M113").Value = -823.4000000000001
$ws.Range("N113").Value = -7330
$ws.Range("H132").Value = 3476.4614
$ws.Range("I132").Value = 3396.75
$ws.Range("J132").Value = 3604
$ws.Range("K132").Value = 10190.25
$ws.Range("L132").Value = 10812
$ws.Range("M132").Value = -7660.25
$ws.Range("N132").Value = -15872

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 833
$ws.Range("I16").Value = 833
$ws.Range("K16").Value = 833
$ws.Range("M16").Value = -663
$ws.Range("H22").Value = 2820.8462
$ws.Range("J22").Value = 3727.4
$ws.Range("L22").Value = 3727.4
$ws.Range("N22").Value = -4317.4
$ws.Range("H27").Value = 2820.8462
$ws.Range("J27").Value = 3727.4
$ws.Range("L27").Value = 3727.4
$ws.Range("N27").Value = -3941.4
$ws.Range("H61").Value = 2337.2307
$ws.Range("I61").Value = 2229.1
$ws.Range("K61").Value = 2229.1
$ws.Range("M61").Value = -2027.1
$ws.Range("H93").Value = 1375.2916
$ws.Range("I93").Value = 1404.15
$ws.Range("J93").Value = 1231
$ws.Range("K93").Value = 1404.15
$ws.Range("L93").Value = 1231
$ws.Range("M93").Value = -156.1500000000001
$ws.Range("N93").Value = -3727
$ws.Range("H106").Value = 17077.8
$ws.Range("J106").Value = 17077.8
$ws.Range("L106").Value = 17077.8
$ws.Range("N106").Value = -19601.8
$ws.Range("H113").Value = 2337.2307
$ws.Range("I113").Value = 2229.1
$ws.Range("K113").Value = 2229.1
$ws.Range("M113").Value = -59.09999999999991
$ws.Range("H140").Value = 92809.664
$ws.Range("J140").Value = 92809.664
$ws.Range("L140").Value = 92809.664
$ws.Range("N140").Value = -103169.664

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4874.8
$ws.Range("I96").Value = 4692.1665
$ws.Range("J96").Value = 5148.75
$ws.Range("K96").Value = 4692.1665
$ws.Range("L96").Value = 5148.75
$ws.Range("M96").Value = -3319.1665
$ws.Range("N96").Value = -7894.75
